$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table with style (border/bold/center/top) for the 4 new rows (23-26)
# by copying the existing formatted cell A2 and pasting formats only.
$ws.Range("A2").Copy()
for ($r = 23; $r -le 26; $r++) {
    $ws.Range("A$r").PasteSpecial(-4122)
}

# Update every data row (2-26) with the refreshed model names and metrics
$ws.Range("A2").Value = "model_5_9_24"
$ws.Range("B2").Value = 0.6181525318771928
$ws.Range("C2").Value = -4.707673193309869
$ws.Range("D2").Value = 0.3776375017972801
$ws.Range("E2").Value = -0.1963937273187066
$ws.Range("F2").Value = 0.4225926399230957
$ws.Range("G2").Value = 1.311975717544556
$ws.Range("H2").Value = 1.042464017868042
$ws.Range("I2").Value = 1.185147643089294

$ws.Range("A3").Value = "model_5_9_23"
$ws.Range("B3").Value = 0.6189279126779217
$ws.Range("C3").Value = -4.629919564925045
$ws.Range("D3").Value = 0.3803210859712595
$ws.Range("E3").Value = -0.1847052757925227
$ws.Range("F3").Value = 0.4217345416545868
$ws.Range("G3").Value = 1.294103145599365
$ws.Range("H3").Value = 1.037968993186951
$ws.Range("I3").Value = 1.173569202423096

$ws.Range("A4").Value = "model_5_9_22"
$ws.Range("B4").Value = 0.6201205155580951
$ws.Range("C4").Value = -4.544298155734875
$ws.Range("D4").Value = 0.3842117324190132
$ws.Range("E4").Value = -0.1710910540414248
$ws.Range("F4").Value = 0.4204146564006805
$ws.Range("G4").Value = 1.27442193031311
$ws.Range("H4").Value = 1.031452059745789
$ws.Range("I4").Value = 1.160082817077637

$ws.Range("A5").Value = "model_5_9_21"
$ws.Range("B5").Value = 0.6208575877610054
$ws.Range("C5").Value = -4.457931172713807
$ws.Range("D5").Value = 0.3860738071938522
$ws.Range("E5").Value = -0.158999161492279
$ws.Range("F5").Value = 0.4195989668369293
$ws.Range("G5").Value = 1.254569530487061
$ws.Range("H5").Value = 1.028333187103271
$ws.Range("I5").Value = 1.148104667663574

$ws.Range("A6").Value = "model_5_9_20"
$ws.Range("B6").Value = 0.6215846293025467
$ws.Range("C6").Value = -4.367092309443401
$ws.Range("D6").Value = 0.3877049781151057
$ws.Range("E6").Value = -0.1465429063307215
$ws.Range("F6").Value = 0.4187943339347839
$ws.Range("G6").Value = 1.233689069747925
$ws.Range("H6").Value = 1.025600910186768
$ws.Range("I6").Value = 1.135765552520752

$ws.Range("A7").Value = "model_5_9_19"
$ws.Range("B7").Value = 0.6228084098093926
$ws.Range("C7").Value = -4.265857545805541
$ws.Range("D7").Value = 0.3905074824563063
$ws.Range("E7").Value = -0.1318758940660041
$ws.Range("F7").Value = 0.4174399375915527
$ws.Range("G7").Value = 1.210419178009033
$ws.Range("H7").Value = 1.020906686782837
$ws.Range("I7").Value = 1.121236324310303

$ws.Range("A8").Value = "model_5_9_18"
$ws.Range("B8").Value = 0.6239881935253491
$ws.Range("C8").Value = -4.155682255841429
$ws.Range("D8").Value = 0.3922619739397935
$ws.Range("E8").Value = -0.1169449750922218
$ws.Range("F8").Value = 0.4161342680454254
$ws.Range("G8").Value = 1.185094237327576
$ws.Range("H8").Value = 1.017967820167542
$ws.Range("I8").Value = 1.106445789337158

$ws.Range("A9").Value = "model_5_9_17"
$ws.Range("B9").Value = 0.6253122203190474
$ws.Range("C9").Value = -4.035947157501942
$ws.Range("D9").Value = 0.3939892785244731
$ws.Range("E9").Value = -0.1008623095829866
$ws.Range("F9").Value = 0.4146689772605896
$ws.Range("G9").Value = 1.15757155418396
$ws.Range("H9").Value = 1.015074610710144
$ws.Range("I9").Value = 1.090514302253723

$ws.Range("A10").Value = "model_5_9_16"
$ws.Range("B10").Value = 0.6283297361353466
$ws.Range("C10").Value = -3.898645890289604
$ws.Range("D10").Value = 0.4014821888398842
$ws.Range("E10").Value = -0.07803260499805842
$ws.Range("F10").Value = 0.4113294780254364
$ws.Range("G10").Value = 1.126011252403259
$ws.Range("H10").Value = 1.002523899078369
$ws.Range("I10").Value = 1.067899107933044

$ws.Range("A11").Value = "model_5_9_15"
$ws.Range("B11").Value = 0.6320932925661306
$ws.Range("C11").Value = -3.751704726516362
$ws.Range("D11").Value = 0.4118000086548904
$ws.Range("E11").Value = -0.05177161569969213
$ws.Range("F11").Value = 0.4071643352508545
$ws.Range("G11").Value = 1.092235088348389
$ws.Range("H11").Value = 0.9852414131164551
$ws.Range("I11").Value = 1.041885018348694

$ws.Range("A12").Value = "model_5_9_14"
$ws.Range("B12").Value = 0.638335366191108
$ws.Range("C12").Value = -3.573892153375944
$ws.Range("D12").Value = 0.4296240181533025
$ws.Range("E12").Value = -0.01574516224492784
$ws.Range("F12").Value = 0.4002561569213867
$ws.Range("G12").Value = 1.051362752914429
$ws.Range("H12").Value = 0.9553859829902649
$ws.Range("I12").Value = 1.006197333335876

$ws.Range("A13").Value = "model_5_9_13"
$ws.Range("B13").Value = 0.6471357701459923
$ws.Range("C13").Value = -3.374175166295248
$ws.Range("D13").Value = 0.4568947822007512
$ws.Range("E13").Value = 0.03048832437902493
$ws.Range("F13").Value = 0.3905167281627655
$ws.Range("G13").Value = 1.005455493927002
$ws.Range("H13").Value = 0.9097071886062622
$ws.Range("I13").Value = 0.960398256778717

$ws.Range("A14").Value = "model_5_9_12"
$ws.Range("B14").Value = 0.6544263503905435
$ws.Range("C14").Value = -3.215544795246121
$ws.Range("D14").Value = 0.4832069141820742
$ws.Range("E14").Value = 0.07091324377573038
$ws.Range("F14").Value = 0.3824481964111328
$ws.Range("G14").Value = 0.9689924716949463
$ws.Range("H14").Value = 0.865634024143219
$ws.Range("I14").Value = 0.9203534126281738

$ws.Range("A15").Value = "model_5_9_11"
$ws.Range("B15").Value = 0.6609143584849981
$ws.Range("C15").Value = -3.066524068595105
$ws.Range("D15").Value = 0.5074368752517073
$ws.Range("E15").Value = 0.1085009895054837
$ws.Range("F15").Value = 0.3752679228782654
$ws.Range("G15").Value = 0.9347382783889771
$ws.Range("H15").Value = 0.8250486850738525
$ws.Range("I15").Value = 0.8831189870834351

$ws.Range("A16").Value = "model_5_9_10"
$ws.Range("B16").Value = 0.6681598192599781
$ws.Range("C16").Value = -2.914374962773083
$ws.Range("D16").Value = 0.5359612556028938
$ws.Range("E16").Value = 0.1498885574086387
$ws.Range("F16").Value = 0.3672492802143097
$ws.Range("G16").Value = 0.8997650146484375
$ws.Range("H16").Value = 0.7772699594497681
$ws.Range("I16").Value = 0.8421204686164856

$ws.Range("A17").Value = "model_5_9_9"
$ws.Range("B17").Value = 0.6793067376250219
$ws.Range("C17").Value = -2.732631634945567
$ws.Range("D17").Value = 0.5783854931028161
$ws.Range("E17").Value = 0.205972971948199
$ws.Range("F17").Value = 0.3549129664897919
$ws.Range("G17").Value = 0.8579891324043274
$ws.Range("H17").Value = 0.7062089443206787
$ws.Range("I17").Value = 0.7865632772445679

$ws.Range("A18").Value = "model_5_9_8"
$ws.Range("B18").Value = 0.6854857481613708
$ws.Range("C18").Value = -2.595010859464164
$ws.Range("D18").Value = 0.6057477340201278
$ws.Range("E18").Value = 0.2446524345261051
$ws.Range("F18").Value = 0.3480746150016785
$ws.Range("G18").Value = 0.8263554573059082
$ws.Range("H18").Value = 0.6603768467903137
$ws.Range("I18").Value = 0.7482473254203796

$ws.Range("A19").Value = "model_5_9_7"
$ws.Range("B19").Value = 0.6922605929887797
$ws.Range("C19").Value = -2.435294867594889
$ws.Range("D19").Value = 0.63272736319576
$ws.Range("E19").Value = 0.2857393587656289
$ws.Range("F19").Value = 0.3405768573284149
$ws.Range("G19").Value = 0.7896428108215332
$ws.Range("H19").Value = 0.6151857376098633
$ws.Range("I19").Value = 0.7075466513633728

$ws.Range("A20").Value = "model_5_9_6"
$ws.Range("B20").Value = 0.6935916782469502
$ws.Range("C20").Value = -2.307598096961833
$ws.Range("D20").Value = 0.6397562211002592
$ws.Range("E20").Value = 0.3070199746122504
$ws.Range("F20").Value = 0.3391036987304688
$ws.Range("G20").Value = 0.760290265083313
$ws.Range("H20").Value = 0.6034122705459595
$ws.Range("I20").Value = 0.6864660382270813

$ws.Range("A21").Value = "model_5_9_5"
$ws.Range("B21").Value = 0.6947392582193285
$ws.Range("C21").Value = -2.139067238263802
$ws.Range("D21").Value = 0.6395543866097374
$ws.Range("E21").Value = 0.3275632333256223
$ws.Range("F21").Value = 0.3378337025642395
$ws.Range("G21").Value = 0.7215513586997986
$ws.Range("H21").Value = 0.6037503480911255
$ws.Range("I21").Value = 0.666115939617157

$ws.Range("A22").Value = "model_5_9_4"
$ws.Range("B22").Value = 0.7008485929291517
$ws.Range("C22").Value = -1.979023630143834
$ws.Range("D22").Value = 0.661566269968242
$ws.Range("E22").Value = 0.3647391893684241
$ws.Range("F22").Value = 0.3310724794864655
$ws.Range("G22").Value = 0.684763491153717
$ws.Range("H22").Value = 0.5668802261352539
$ws.Range("I22").Value = 0.6292893886566162

$ws.Range("A23").Value = "model_5_9_3"
$ws.Range("B23").Value = 0.7089380036509947
$ws.Range("C23").Value = -1.679883220482193
$ws.Range("D23").Value = 0.6721401587725016
$ws.Range("E23").Value = 0.409900542792368
$ws.Range("F23").Value = 0.3221198618412018
$ws.Range("G23").Value = 0.61600261926651
$ws.Range("H23").Value = 0.5491688251495361
$ws.Range("I23").Value = 0.584552526473999

$ws.Range("A24").Value = "model_5_9_2"
$ws.Range("B24").Value = 0.7182267741326472
$ws.Range("C24").Value = -1.644677225653556
$ws.Range("D24").Value = 0.7384073191376359
$ws.Range("E24").Value = 0.4669571406973976
$ws.Range("F24").Value = 0.3118399083614349
$ws.Range("G24").Value = 0.6079100370407104
$ws.Range("H24").Value = 0.4381706118583679
$ws.Range("I24").Value = 0.5280323028564453

$ws.Range("A25").Value = "model_5_9_0"
$ws.Range("B25").Value = 0.7796108279776712
$ws.Range("C25").Value = 0.8996436517196014
$ws.Range("D25").Value = 0.925077900726047
$ws.Range("E25").Value = 0.9280545078748007
$ws.Range("F25").Value = 0.2439058572053909
$ws.Range("G25").Value = 0.02306808531284332
$ws.Range("H25").Value = 0.1254953444004059
$ws.Range("I25").Value = 0.07126921415328979

$ws.Range("A26").Value = "model_5_9_1"
$ws.Range("B26").Value = 0.7855654854413647
$ws.Range("C26").Value = 0.1894413810049099
$ws.Range("D26").Value = 0.89336157548061
$ws.Range("E26").Value = 0.8155710763499575
$ws.Range("F26").Value = 0.2373158186674118
$ws.Range("G26").Value = 0.1863164156675339
$ws.Range("H26").Value = 0.1786205321550369
$ws.Range("I26").Value = 0.1826952993869781
